# Swap the "Enterprises (absolute #)" and "Enterprises density (per 1000 people)"
# rows in both indicator tables on the sheet, so that the "density" row now
# appears before the "absolute number" row (matching the row order used by
# the target workbook revision).
#
# Note: the numeric-looking values in these rows (e.g. "429035", "20.1") are
# stored as TEXT in the workbook, not as numbers. A plain ".Value = ..."
# assignment of such a string would be auto-coerced to a number by Excel,
# changing the cell's stored type. To preserve the original text storage
# exactly, cell contents are exchanged using Copy / PasteSpecial (values
# only) through a scratch cell, which keeps the text type intact and does
# not introduce any new cell styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

function Swap-CellValues($ws, $ref1, $ref2, $scratch) {
    $c1 = $ws.Range($ref1)
    $c2 = $ws.Range($ref2)

    $c1.Copy()
    $scratch.PasteSpecial($xlPasteValues)

    $c2.Copy()
    $c1.PasteSpecial($xlPasteValues)

    $scratch.Copy()
    $c2.PasteSpecial($xlPasteValues)

    $scratch.ClearContents()
}

# Scratch cell far away from any real data, used as temporary holding spot.
$scratch = $ws.Range("Z100")

# --- Table 1: "Source Type: Statistical Institution" (rows 10-11) ---
Swap-CellValues $ws "A10" "A11" $scratch
Swap-CellValues $ws "D10" "D11" $scratch

# --- Table 2: "Source Type: SME Associations (Most Widely Used)" (rows 32-33) ---
Swap-CellValues $ws "A32" "A33" $scratch
Swap-CellValues $ws "B32" "B33" $scratch
Swap-CellValues $ws "C32" "C33" $scratch
Swap-CellValues $ws "D32" "D33" $scratch
